$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.311.06"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.589.97"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "209.86"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").Value = "19.52"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.586.62"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "64.45"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "26.322.39"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("E19").Value = "  +4.53%  "
$ws.Range("D20").Value = "211.02"
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "2.15"
$ws.Range("E23").Value = "  -4.76%  "
$ws.Range("D24").Value = "8.92"
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("D25").Value = "144.83"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "15.28"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "1.306.10"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("E35").Value = "  +2.94%  "
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("E37").Value = "  -0.68%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  -13.75%  "
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("E42").Value = "  +3.56%  "
$ws.Range("D43").Value = "62.58"
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.13"
$ws.Range("E44").Value = "  -4.23%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "0.764"
$ws.Range("E45").Value = "  -1.78%  "
$ws.Range("D46").Value = "1.725.92"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").Value = "87.97"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("E48").Value = "  -4.35%  "
$ws.Range("D49").Value = "0.0981"
$ws.Range("E49").Value = "  -4.38%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("E51").Value = "  -0.52%  "
